$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Tff3"
$ws.Range("C2").Value = "Ackr3"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 0.7627916666666668
$ws.Range("H2").Value = 2.288375
$ws.Range("I2").Value = 0.6223788491141802
$ws.Range("J2").Value = 0.6223788491141802
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 15.65740233333333
$ws.Range("N2").Value = 46.972207
$ws.Range("O2").Value = 0.5111560914107862
$ws.Range("P2").Value = 0.5111560914107862
$ws.Range("Q2").Value = 11.94333602151389
$ws.Range("R2").Value = 107.490024193625
$ws.Range("S2").Value = 0.3181327398899478
$ws.Range("T2").Value = 0.3181327398899478

$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Tff3"
$ws.Range("C3").Value = "Ackr3"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 0.7627916666666668
$ws.Range("H3").Value = 2.288375
$ws.Range("I3").Value = 0.6223788491141802
$ws.Range("J3").Value = 0.6223788491141802
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 8.444501666666666
$ws.Range("N3").Value = 25.333505
$ws.Range("O3").Value = 0.275681647182037
$ws.Range("P3").Value = 0.275681647182037
$ws.Range("Q3").Value = 6.441395500486111
$ws.Range("R3").Value = 57.97255950437501
$ws.Range("S3").Value = 0.1715784262950577
$ws.Range("T3").Value = 0.1715784262950577

$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Tff3"
$ws.Range("C4").Value = "Ackr3"
$ws.Range("D4").Value = "M1"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 0.7627916666666668
$ws.Range("H4").Value = 2.288375
$ws.Range("I4").Value = 0.6223788491141802
$ws.Range("J4").Value = 0.6223788491141802
$ws.Range("K4").Value = 2
$ws.Range("L4").Value = 0.6666666666666666
$ws.Range("M4").Value = 0.186037
$ws.Range("N4").Value = 0.558111
$ws.Range("O4").Value = 0.006073417783698461
$ws.Range("P4").Value = 0.006073417783698461
$ws.Range("Q4").Value = 0.1419074732916667
$ws.Range("R4").Value = 1.277167259625
$ws.Range("S4").Value = 0.003779966770407843
$ws.Range("T4").Value = 0.003779966770407843

$ws.Range("A5").Value = "ECs"
$ws.Range("B5").Value = "Tff3"
$ws.Range("C5").Value = "Ackr3"
$ws.Range("D5").Value = "sCs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 0.7627916666666668
$ws.Range("H5").Value = 2.288375
$ws.Range("I5").Value = 0.6223788491141802
$ws.Range("J5").Value = 0.6223788491141802
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 6.343411333333333
$ws.Range("N5").Value = 19.030234
$ws.Range("O5").Value = 0.2070888436234783
$ws.Range("P5").Value = 0.2070888436234783
$ws.Range("Q5").Value = 4.838701303305556
$ws.Range("R5").Value = 43.54831172975
$ws.Range("S5").Value = 0.1288877161587669
$ws.Range("T5").Value = 0.1288877161587668

$ws.Range("A6").Value = "M2"
$ws.Range("B6").Value = "Tff3"
$ws.Range("C6").Value = "Ackr3"
$ws.Range("D6").Value = "ECs"
$ws.Range("E6").Value = 1
$ws.Range("F6").Value = 0.3333333333333333
$ws.Range("G6").Value = 0.462815
$ws.Range("H6").Value = 1.388445
$ws.Range("I6").Value = 0.3776211508858198
$ws.Range("J6").Value = 0.3776211508858198
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 15.65740233333333
$ws.Range("N6").Value = 46.972207
$ws.Range("O6").Value = 0.5111560914107862
$ws.Range("P6").Value = 0.5111560914107862
$ws.Range("Q6").Value = 7.246480660901666
$ws.Range("R6").Value = 65.218325948115
$ws.Range("S6").Value = 0.1930233515208384
$ws.Range("T6").Value = 0.1930233515208384

$ws.Range("A7").Value = "M2"
$ws.Range("B7").Value = "Tff3"
$ws.Range("C7").Value = "Ackr3"
$ws.Range("D7").Value = "FAPs"
$ws.Range("E7").Value = 1
$ws.Range("F7").Value = 0.3333333333333333
$ws.Range("G7").Value = 0.462815
$ws.Range("H7").Value = 1.388445
$ws.Range("I7").Value = 0.3776211508858198
$ws.Range("J7").Value = 0.3776211508858198
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 8.444501666666666
$ws.Range("N7").Value = 25.333505
$ws.Range("O7").Value = 0.275681647182037
$ws.Range("P7").Value = 0.275681647182037
$ws.Range("Q7").Value = 3.908242038858333
$ws.Range("R7").Value = 35.174178349725
$ws.Range("S7").Value = 0.1041032208869793
$ws.Range("T7").Value = 0.1041032208869793

$ws.Range("A8").Value = "M2"
$ws.Range("B8").Value = "Tff3"
$ws.Range("C8").Value = "Ackr3"
$ws.Range("D8").Value = "M1"
$ws.Range("E8").Value = 1
$ws.Range("F8").Value = 0.3333333333333333
$ws.Range("G8").Value = 0.462815
$ws.Range("H8").Value = 1.388445
$ws.Range("I8").Value = 0.3776211508858198
$ws.Range("J8").Value = 0.3776211508858198
$ws.Range("K8").Value = 2
$ws.Range("L8").Value = 0.6666666666666666
$ws.Range("M8").Value = 0.186037
$ws.Range("N8").Value = 0.558111
$ws.Range("O8").Value = 0.006073417783698461
$ws.Range("P8").Value = 0.006073417783698461
$ws.Range("Q8").Value = 0.086100714155
$ws.Range("R8").Value = 0.774906427395
$ws.Range("S8").Value = 0.002293451013290618
$ws.Range("T8").Value = 0.002293451013290617

$ws.Range("A9").Value = "M2"
$ws.Range("B9").Value = "Tff3"
$ws.Range("C9").Value = "Ackr3"
$ws.Range("D9").Value = "sCs"
$ws.Range("E9").Value = 1
$ws.Range("F9").Value = 0.3333333333333333
$ws.Range("G9").Value = 0.462815
$ws.Range("H9").Value = 1.388445
$ws.Range("I9").Value = 0.3776211508858198
$ws.Range("J9").Value = 0.3776211508858198
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 6.343411333333333
$ws.Range("N9").Value = 19.030234
$ws.Range("O9").Value = 0.2070888436234783
$ws.Range("P9").Value = 0.2070888436234783
$ws.Range("Q9").Value = 2.935825916236666
$ws.Range("R9").Value = 26.42243324613
$ws.Range("S9").Value = 0.07820112746471145
$ws.Range("T9").Value = 0.07820112746471142
